$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.48670494556427
$ws.Range("B1").Value = 1.646387457847595
$ws.Range("C1").Value = 1.718206882476807
$ws.Range("D1").Value = 2.308779954910278
$ws.Range("E1").Value = 3.798918724060059
